$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "4 KBs, 24 query sets." -> "3 KBS, 15 Query sets."
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "4 KBs, 24 query sets.", $true, $false, $false, $false, $false,
    $true, 1, $false, "3 KBS, 15 Query sets.", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "For each number of ranks (3, 10, 50, 100), " ->
#    "For each number of ranks (10, 50, 100), "
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "For each number of ranks (3, 10, 50, 100), ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "For each number of ranks (10, 50, 100), ", 2) | Out-Null

# ------------------------------------------------------------------
# 3) " distribution, same number of ranks in kb (100))" ->
#    " distribution, same number of ranks in kb (" / "5" / "0))"
#    (the "10" inside "100" is replaced by "5", leaving the run split
#    into three pieces the way Word splits a run when you select and
#    retype only part of its text)
# ------------------------------------------------------------------
$found = $d.Content.Find
$found.Execute("distribution, same number of ranks in kb (100))") | Out-Null
$matchRange = $found.Parent
$matchStart = $matchRange.Start
$offset = $matchRange.Text.IndexOf("100")
$tenStart = $matchStart + $offset
$tenEnd = $tenStart + 2
$tenRange = $d.Range($tenStart, $tenEnd)
$tenRange.Bold = 1
$tenRange.Text = "5"
$tenRange.Bold = 0

# ------------------------------------------------------------------
# 4) "(Same number of ranks in kb (100))" ->
#    "(Same number of ranks in kb (" / "5" / "0))"
# ------------------------------------------------------------------
$found2 = $d.Content.Find
$found2.Execute("(Same number of ranks in kb (100))") | Out-Null
$matchRange2 = $found2.Parent
$matchStart2 = $matchRange2.Start
$offset2 = $matchRange2.Text.IndexOf("100")
$tenStart2 = $matchStart2 + $offset2
$tenEnd2 = $tenStart2 + 2
$tenRange2 = $d.Range($tenStart2, $tenEnd2)
$tenRange2.Bold = 1
$tenRange2.Text = "5"
$tenRange2.Bold = 0
